# emails_muq.xlsx edit
# - B2 email changes to anantdeep.parihar@mu-sigma.com and becomes a mailto: hyperlink
# - C2 email changes to cs.satish@mu-sigma.com (new value) and becomes a mailto: hyperlink
# - B3 keeps its existing email (kumar.singh@mu-sigma.com) but becomes a mailto: hyperlink
# - Columns A and B are widened
# - Selection moves to C2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Hyperlinks / email value updates (order matters: it drives the r:id numbering) ---

# C2: new email address, turned into a mailto hyperlink
$ws.Range("C2").Value = "cs.satish@mu-sigma.com"
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:cs.satish@mu-sigma.com")

# B3: same email as before, now wired up as a mailto hyperlink
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:kumar.singh@mu-sigma.com")

# B2: new email address, turned into a mailto hyperlink
$ws.Range("B2").Value = "anantdeep.parihar@mu-sigma.com"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:anantdeep.parihar@mu-sigma.com")

# --- Column widths ---
$ws.Columns("A").ColumnWidth = 30.666666666666668
$ws.Columns("B").ColumnWidth = 32.166666666666664

# --- Selection ---
$ws.Range("C2").Select()
